$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.811.82"
$ws.Range("E2").Value = "  +4.65%  "

# Row 3
$ws.Range("D3").Value = "2.276.01"
$ws.Range("E3").Value = "  +2.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6
$ws.Range("E6").Value = "  +0.70%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.420"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.90%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0948"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.75%  "

# Row 11
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "

# Row 12
$ws.Range("E12").Value = "  +0.77%  "

# Row 13
$ws.Range("D13").Value = "2.614.24"
$ws.Range("E13").Value = "  +2.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18
$ws.Range("D18").Value = "2.280.97"
$ws.Range("E18").Value = "  +2.40%  "

# Row 19
$ws.Range("D19").Value = "43.714.64"
$ws.Range("E19").Value = "  +4.73%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  +5.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("E25").Value = "  +7.12%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.27%  "

# Row 29
$ws.Range("E29").Value = "  -0.35%  "

# Row 30
$ws.Range("E30").Value = "  +3.36%  "

# Row 31
$ws.Range("E31").Value = "  +5.56%  "

# Row 32
$ws.Range("E32").Value = "  +0.77%  "

# Row 33
$ws.Range("E33").Value = "  +0.47%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.45%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0661"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.71%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "

# Row 38
$ws.Range("E38").Value = "  +2.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0251"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.89%  "

# Row 41
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000225"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -13.10%  "

# Row 44
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.51%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0970"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "

# Row 48
$ws.Range("D48").Value = "1.470.67"
$ws.Range("E48").Value = "  +0.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "

# Row 50
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.31%  "

# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
